# Applies the "Add files via upload" revision to the Absorber Section
# Heat Duty and Tray Counts workbook:
#   - Renames the compressor block in B1 from "C-601B" to "C-601A"
#   - Adds a single horsepower reading in B12
#   - Strips the bracketed units off several row labels in column A
#     (e.g. "Indicated horsepower [kW]" -> "Indicated horsepower")
#   - Clears the now-stale numeric readouts that used to sit in column B
#     next to those relabeled rows
#   - Moves the active selection to B19
#   - Restores the (smaller) saved window size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: compressor name -------------------------------------------------
$ws.Range("B1").Value = "C-601A"

# --- B12: new reading ------------------------------------------------------
$ws.Range("B12").Value = 84

# --- Column A relabels (drop the trailing "[unit]") + clear old B values ---
$ws.Range("A17").Value = "Indicated horsepower"
$ws.Range("B17").Value = $null

$ws.Range("A18").Value = "Calculated brake horsepower"
$ws.Range("B18").Value = $null

$ws.Range("A19").Value = "Net work required"
$ws.Range("B19").Value = $null

$ws.Range("A20").Value = "Power loss"
$ws.Range("B20").Value = $null

# A21 keeps its existing label ("Efficiency (polytropic / isentropic) used")
$ws.Range("B21").Value = $null

$ws.Range("A22").Value = "Calculated discharge pressure"
$ws.Range("B22").Value = $null

$ws.Range("A23").Value = "Calculated pressure change"
$ws.Range("B23").Value = $null

# A24 keeps its existing label ("Calculated pressure ratio")
$ws.Range("B24").Value = $null

$ws.Range("A25").Value = "Outlet temperature"
$ws.Range("B25").Value = $null

$ws.Range("A26").Value = "Isentropic outlet temperature"
$ws.Range("B26").Value = $null

# A27 keeps its existing label ("Vapor fraction")
$ws.Range("B27").Value = $null

$ws.Range("A30").Value = "Head developed"
$ws.Range("B30").Value = $null

$ws.Range("A31").Value = "Isentropic power requirement"
$ws.Range("B31").Value = $null

# A32 keeps its existing label ("Inlet heat capacity ratio")
$ws.Range("B32").Value = $null

$ws.Range("A33").Value = "Inlet volumetric flow rate"
$ws.Range("B33").Value = $null

$ws.Range("A34").Value = "Outlet volumetric flow rate"
$ws.Range("B34").Value = $null

# A35 / A36 keep their existing labels (compressibility factors)
$ws.Range("B35").Value = $null
$ws.Range("B36").Value = $null

$ws.Range("A44").Value = "Total feed stream CO2e flow"
$ws.Range("B44").Value = $null

$ws.Range("A45").Value = "Total product stream CO2e flow"
$ws.Range("B45").Value = $null

$ws.Range("A46").Value = "Net stream CO2e production"
$ws.Range("B46").Value = $null

$ws.Range("A47").Value = "Utility CO2e production"
$ws.Range("B47").Value = $null

$ws.Range("A48").Value = "Total CO2e production"
$ws.Range("B48").Value = $null

# --- Selection + saved window geometry -------------------------------------
$ws.Range("B19").Select() | Out-Null

# Best-effort: mirror the saved window size from the workbook's bookViews
# (xWindow/yWindow stay 0; width/height shrink to 13478 x 12608). Some hosts
# don't expose independent window geometry, so this is harmless if ignored.
$excel.ActiveWindow.Width = 13478
$excel.ActiveWindow.Height = 12608
